$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.691262245178223
$ws.Range("B1").Value = 5.284496307373047
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 4.97516918182373
$ws.Range("E1").Value = 2.409585952758789
